# Apply newly added iAuthor TC's credential updates to the "users" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds a single candidate's generated credentials; update them
# to the freshly issued values while leaving Title (E2) and Role (H2)
# untouched.
$ws.Range("A2").Value = "rNNJp810"   # Client Id
$ws.Range("B2").Value = 231102296    # Candidate ID
$ws.Range("C2").Value = "qgfjyfj84"  # User Name
$ws.Range("D2").Value = "xJQ6&%5s"   # Exam Password
$ws.Range("F2").Value = "GezVubln"   # First Name
$ws.Range("G2").Value = "QDYm"       # Last Name
